$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) / volume-change (E) cells with the latest scraped snapshot values.
# D-column numeric-looking strings (e.g. "6.70", "0.0000138") are entered with a leading
# apostrophe so Excel keeps them as text (matching the source data's inlineStr/General
# cells) instead of auto-converting them to numbers and dropping formatting like trailing zeros.

$ws.Range('D2').Value = '59.729.37'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '2.599.82'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'556.93"
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').Value = "'141.64"
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('D7').Value = "'0.997"
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = "'0.598"
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = '2.620.66'
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('D10').Value = "'6.70"
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('D11').Value = "'0.105"
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').Value = "'0.162"
$ws.Range('E12').Value = '  +6.73%  '
$ws.Range('E13').Value = '  +7.45%  '
$ws.Range('D14').Value = '3.058.00'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').Value = "'23.48"
$ws.Range('E15').Value = '  +6.44%  '
$ws.Range('D16').Value = '59.703.49'
$ws.Range('E16').Value = '  +1.50%  '
$ws.Range('D17').Value = "'0.0000138"
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').Value = '2.606.42'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('D19').Value = "'4.60"
$ws.Range('E19').Value = '  +2.69%  '
$ws.Range('D20').Value = "'342.23"
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('D21').Value = "'10.62"
$ws.Range('E21').Value = '  +4.68%  '
$ws.Range('D22').Value = "'6.73"
$ws.Range('E22').Value = '  +9.60%  '
$ws.Range('D23').Value = "'0.999"
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = "'0.518"
$ws.Range('E24').Value = '  +14.79%  '
$ws.Range('D25').Value = "'62.40"
$ws.Range('E25').Value = '  -2.26%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('D28').Value = "'7.57"
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('D29').Value = '0.0₃0782'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('D30').Value = "'0.997"
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('D32').Value = "'6.18"
$ws.Range('E32').Value = '  +2.62%  '
$ws.Range('D33').Value = "'159.01"
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('D34').Value = "'19.38"
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('D35').Value = "'4.12"
$ws.Range('E35').Value = '  +3.39%  '
$ws.Range('D36').Value = "'0.919"
$ws.Range('E36').Value = '  +5.03%  '
$ws.Range('D37').Value = "'1.18"
$ws.Range('E37').Value = '  +4.91%  '
$ws.Range('E38').Value = '  +2.84%  '
$ws.Range('E39').Value = '  +2.27%  '
$ws.Range('D40').Value = "'0.846"
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('D41').Value = "'3.70"
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('D42').Value = "'291.41"
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').Value = "'139.26"
$ws.Range('E43').Value = '  +12.54%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = "'0.0980"
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('D46').Value = "'0.603"
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('D47').Value = "'0.0240"
$ws.Range('E47').Value = '  +3.85%  '
$ws.Range('D48').Value = "'0.0539"
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('E49').Value = '  +0.21%  '

# Row 50/51 shift: RenderToken inserted as new #48, InjectiveProtocol moves down to #49
# (price/volume refreshed), Maker (former #49) drops off the list.
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'4.80"
$ws.Range('E50').Value = '  +7.37%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = "'18.99"
$ws.Range('E51').Value = '  +2.82%  '
